$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the per-vendor "Particulars" labels in column A with their
# chemical-category grouping. Column B (Credit amounts) is untouched.
$categories = @{
    2  = "Polymers"
    3  = "Basic Chemicals"
    4  = "Cosmetic Chemicals"
    5  = "Cosmetic Chemicals"
    6  = "Basic Chemicals"
    7  = "Cosmetic Chemicals"
    8  = "Specialty Chemicals"
    9  = "Pharmaceuticals"
    10 = "Cosmetic Chemicals"
    11 = "Pesticides"
    12 = "Pharmaceuticals"
    13 = "Pharmaceuticals"
    14 = "Pesticides"
    15 = "Biochemicals"
    16 = "Pharmaceuticals"
    17 = "Pesticides"
    18 = "Polymers"
    19 = "Pesticides"
}

foreach ($row in $categories.Keys) {
    $ws.Cells.Item($row, 1).Value = $categories[$row]
}
